# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for rows 2-19 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 6
    4  = 6
    5  = 3
    6  = 4
    7  = 0
    8  = 3
    9  = 5
    10 = 9
    11 = 3
    12 = 10
    13 = 3
    14 = 3
    15 = 10
    16 = 5
    17 = 2
    18 = 1
    19 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
